# Generate Report for Handback
# Refresh the handoff/handback timestamps for the 95daaf45-... row
# (the row whose Status is "Handed back: in sync with en-US"):
#   - Overview!G3        Correspond Handback DateTime
#   - zh-cn!H3            Correspond Handoff Datetime
#   - zh-cn!K3            Correspond Handback DateTime
#   - de-de!H3            Correspond Handoff Datetime (same text as Overview!G3)
#   - de-de!K3            Correspond Handback DateTime

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-17 18:44:49"

$wsZhCn.Range("H3").Value = "2016-08-17 18:44:44"
$wsZhCn.Range("K3").Value = "2016-08-17 18:45:08"

# de-de sheet's "Correspond Handoff Datetime" (H3) shares the same text as the
# Overview sheet's "Correspond Handback DateTime" (G3) for this row, so update
# it in lockstep to keep both cells showing the same, refreshed timestamp.
$wsDeDe.Range("H3").Value = "2016-08-17 18:44:49"
$wsDeDe.Range("K3").Value = "2016-08-17 18:45:18"
